$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A19").Value = "18_021215_2150_7_xgboost_binary_logits_with_more_random_3in1_preprocess_valid1_valid2_"
$ws.Range("B19").Value = 0.6312
$ws.Range("C19").Value = "ensembled 7 tree xgboost binary logits on more random combined 3in1 data set with features preprocessed, with 2 valid sets"
$ws.Range("C20").Select()
